# "modal resultados para envio" - add a new "Enviar Datos" (send data) confirmation
# modal's text rows to the Translation sheet, drop a handful of obsolete rows, and
# fix the alignment of the result "<value>" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------------
# 1. Fix alignment of the "<value>" label (widget 342, originally row 22)
#    from Center to Left. Do this before any row deletions so the row index
#    is still the original one.
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = "Left"

# ---------------------------------------------------------------------------
# 2. Remove obsolete rows (processed from the bottom up so row numbers of
#    not-yet-deleted rows are unaffected).
#    Row 46 -> widget 607 "OK"
#    Row 45 -> widget 605 "Se envio el\nregistro con éxito"
#    Row 12 -> widget 186 "Administrador de los perfiles\nPor productos..."
#    Row 10 -> widget 174 "Empezar"
#    Row 9  -> widget 172 "Perfiles"
# ---------------------------------------------------------------------------
$ws.Rows.Item(46).Delete()
$ws.Rows.Item(45).Delete()
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()

# ---------------------------------------------------------------------------
# 3. Append the new "send data" modal rows (these land on rows 56-70 after
#    the deletions above).
# ---------------------------------------------------------------------------
$ws.Range("B56").Value = "SingleUseId140"
$ws.Range("C56").Value = "menu_main"
$ws.Range("D56").Value = "Left"
$ws.Range("E56").Value = "LTR"
$ws.Range("F56").Value = "Enviar Datos"

$ws.Range("B57").Value = "SingleUseId141"
$ws.Range("C57").Value = "Default"
$ws.Range("D57").Value = "Left"
$ws.Range("E57").Value = "LTR"
$ws.Range("F57").Value = "Peso Obtenido"

$ws.Range("B58").Value = "SingleUseId142"
$ws.Range("C58").Value = "perfilesbtn"
$ws.Range("D58").Value = "Left"
$ws.Range("E58").Value = "LTR"
$ws.Range("F58").Value = "peso esperado"

$ws.Range("B59").Value = "SingleUseId143"
$ws.Range("C59").Value = "perfilesbtn"
$ws.Range("D59").Value = "Left"
$ws.Range("E59").Value = "LTR"
$ws.Range("F59").Value = "diferencia obtenida"

$ws.Range("B60").Value = "SingleUseId144"
$ws.Range("C60").Value = "perfilesbtn"
$ws.Range("D60").Value = "Left"
$ws.Range("E60").Value = "LTR"
$ws.Range("F60").Value = "diferencia permitida"

$ws.Range("B61").Value = "SingleUseId145"
$ws.Range("C61").Value = "Default"
$ws.Range("D61").Value = "Left"
$ws.Range("E61").Value = "LTR"
$ws.Range("F61").Value = "<d> [kg]"

$ws.Range("B62").Value = "SingleUseId146"
$ws.Range("C62").Value = "Default"
$ws.Range("D62").Value = "Left"
$ws.Range("E62").Value = "LTR"
$ws.Range("F62").Value = "00.00"

$ws.Range("B63").Value = "SingleUseId147"
$ws.Range("C63").Value = "perfilesbtn"
$ws.Range("D63").Value = "Right"
$ws.Range("E63").Value = "LTR"
$ws.Range("F63").Value = "<d> [kg]"

$ws.Range("B64").Value = "SingleUseId148"
$ws.Range("C64").Value = "perfilesbtn"
$ws.Range("D64").Value = "Left"
$ws.Range("E64").Value = "LTR"
$ws.Range("F64").Value = "00.00"

$ws.Range("B65").Value = "SingleUseId149"
$ws.Range("C65").Value = "perfilesbtn"
$ws.Range("D65").Value = "Right"
$ws.Range("E65").Value = "LTR"
$ws.Range("F65").Value = "<d> [kg]"

$ws.Range("B66").Value = "SingleUseId150"
$ws.Range("C66").Value = "perfilesbtn"
$ws.Range("D66").Value = "Left"
$ws.Range("E66").Value = "LTR"
$ws.Range("F66").Value = "00.00"

$ws.Range("B67").Value = "SingleUseId151"
$ws.Range("C67").Value = "perfilesbtn"
$ws.Range("D67").Value = "Right"
$ws.Range("E67").Value = "LTR"
$ws.Range("F67").Value = "<d> [kg]"

$ws.Range("B68").Value = "SingleUseId152"
$ws.Range("C68").Value = "perfilesbtn"
$ws.Range("D68").Value = "Left"
$ws.Range("E68").Value = "LTR"
$ws.Range("F68").Value = "00.00"

$ws.Range("B69").Value = "SingleUseId153"
$ws.Range("C69").Value = "perfilesbtn"
$ws.Range("D69").Value = "Center"
$ws.Range("E69").Value = "LTR"
$ws.Range("F69").Value = "    Enviar"

$ws.Range("B70").Value = "SingleUseId154"
$ws.Range("C70").Value = "perfilesbtn"
$ws.Range("D70").Value = "Left"
$ws.Range("E70").Value = "LTR"
$ws.Range("F70").Value = "Datos `nInvalidos"
